$wb = $excel.ActiveWorkbook

# Cell value updates per sheet, derived from the authoritative diff
# (Moogle market-data refresh: currentAveragePrice* / LevePrice* / LeveProfit* columns)
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = [ordered]@{
  "H18" = 859
  "I18" = 859
  "K18" = 859
  "M18" = -575
  "H19" = 898.4666999999999
  "I19" = 562.6
  "J19" = 1066.4
  "K19" = 562.6
  "L19" = 1066.4
  "M19" = -387.6
  "N19" = -1416.4
  "H28" = 845.5294
  "J28" = 1815.6
  "L28" = 1815.6
  "N28" = -2785.6
  "H96" = 315.5
  "J96" = 0
  "L96" = 0
  "H107" = 1577.2307
  "I107" = 1349.8
  "K107" = 1349.8
  "M107" = 570.2
  "H112" = 6170
  "I112" = 2444.5
  "J112" = 7101.375
  "K112" = 7333.5
  "L112" = 21304.125
  "M112" = -6225.5
  "N112" = -23520.125
  "H125" = 64764.117
  "J125" = 118266.11
  "L125" = 1064394.99
  "N125" = -1069314.99
}
foreach ($cell in $ALC_updates.Keys) {
  $ws.Range($cell).Value = $ALC_updates[$cell]
}

$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = [ordered]@{
  "H2" = 777.8461
  "I2" = 784.3333
  "J2" = 700
  "K2" = 784.3333
  "L2" = 700
  "M2" = -671.3333
  "N2" = -926
  "H32" = 9143.322
  "I32" = 5654.079
  "K32" = 5654.079
  "M32" = -5367.079
  "H45" = 2895.5715
  "I45" = 2321.818
  "K45" = 2321.818
  "M45" = -1944.818
  "H110" = 1164.1
  "I110" = 1167.625
  "K110" = 1167.625
  "M110" = 877.375
  "H116" = 777.8461
  "I116" = 784.3333
  "J116" = 700
  "K116" = 784.3333
  "L116" = 700
  "M116" = 1509.6667
  "N116" = -5288
}
foreach ($cell in $ARM_updates.Keys) {
  $ws.Range($cell).Value = $ARM_updates[$cell]
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = [ordered]@{
  "H3" = 777.8461
  "I3" = 784.3333
  "J3" = 700
  "K3" = 784.3333
  "L3" = 700
  "M3" = -670.3333
  "N3" = -928
  "H80" = 16941.834
  "I80" = 50025.5
  "J80" = 10325.1
  "K80" = 50025.5
  "L80" = 10325.1
  "M80" = -49027.5
  "N80" = -12321.1
  "H83" = 16941.834
  "I83" = 50025.5
  "J83" = 10325.1
  "K83" = 250127.5
  "L83" = 51625.5
  "M83" = -245135.5
  "N83" = -61609.5
  "H105" = 3795.7188
  "I105" = 4028.44
  "K105" = 4028.44
  "M105" = -2281.44
  "H134" = 4245.8823
  "I134" = 3278.6667
  "K134" = 9836.000100000001
  "M134" = -7301.000100000001
}
foreach ($cell in $BSM_updates.Keys) {
  $ws.Range($cell).Value = $BSM_updates[$cell]
}

$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = [ordered]@{
  "H16" = 1013.0625
  "I16" = 897.1818
  "J16" = 1268
  "K16" = 897.1818
  "L16" = 1268
  "M16" = -610.1818
  "N16" = -1842
  "H58" = 2368.2222
  "I58" = 2209.8667
  "J58" = 2566.1667
  "K58" = 2209.8667
  "L58" = 2566.1667
  "M58" = -2006.8667
  "N58" = -2972.1667
  "H107" = 3649.3572
  "I107" = 3756.5
  "J107" = 3006.5
  "K107" = 3756.5
  "L107" = 3006.5
  "M107" = -1836.5
  "N107" = -6846.5
  "H113" = 1013.0625
  "I113" = 897.1818
  "J113" = 1268
  "K113" = 897.1818
  "L113" = 1268
  "M113" = 1272.8182
  "N113" = -5608
  "H122" = 4127.8
  "I122" = 2988.25
  "K122" = 8964.75
  "M122" = -6514.75
  "H134" = 2297.3928
  "I134" = 1677.7778
  "K134" = 5033.3334
  "M134" = -2498.3334
  "H136" = 2368.2222
  "I136" = 2209.8667
  "J136" = 2566.1667
  "K136" = 6629.6001
  "L136" = 7698.500100000001
  "M136" = -4079.6001
  "N136" = -12798.5001
}
foreach ($cell in $CRP_updates.Keys) {
  $ws.Range($cell).Value = $CRP_updates[$cell]
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = [ordered]@{
  "H131" = 3128.3333
  "I131" = 1605.9231
  "J131" = 3988.8262
  "K131" = 4817.7693
  "L131" = 11966.4786
  "M131" = 222.2307000000001
  "N131" = -22046.4786
  "H139" = 2311.6
  "I139" = 2268.4443
  "K139" = 6805.3329
  "M139" = -1665.3329
}
foreach ($cell in $CUL_updates.Keys) {
  $ws.Range($cell).Value = $CUL_updates[$cell]
}

$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = [ordered]@{
  "H70" = 7499
  "J70" = 7499
  "L70" = 7499
  "N70" = -8039
  "H73" = 7499
  "J73" = 7499
  "L73" = 7499
  "N73" = -9371
  "H80" = 2540.923
  "I80" = 2336
  "K80" = 2336
  "M80" = -1338
  "H83" = 2540.923
  "I83" = 2336
  "K83" = 11680
  "M83" = -6688
  "H102" = 1823.9474
  "I102" = 1154.0667
  "J102" = 4336
  "K102" = 1154.0667
  "L102" = 4336
  "M102" = 467.9332999999999
  "N102" = -7580
}
foreach ($cell in $GSM_updates.Keys) {
  $ws.Range($cell).Value = $GSM_updates[$cell]
}

$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = [ordered]@{
  "H7" = 96952
  "I7" = 129371
  "K7" = 129371
  "M7" = -129259
  "H61" = 6344.9165
  "I61" = 2314.4285
  "J61" = 11987.6
  "K61" = 2314.4285
  "L61" = 11987.6
  "M61" = -2112.4285
  "N61" = -12391.6
  "H113" = 6344.9165
  "I113" = 2314.4285
  "J113" = 11987.6
  "K113" = 2314.4285
  "L113" = 11987.6
  "M113" = -144.4285
  "N113" = -16327.6
  "H126" = 96952
  "I126" = 129371
  "K126" = 388113
  "M126" = -385643
}
foreach ($cell in $LTW_updates.Keys) {
  $ws.Range($cell).Value = $LTW_updates[$cell]
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = [ordered]@{
  "H122" = 3024.2856
  "I122" = 2500.5789
  "K122" = 7501.736699999999
  "M122" = -5051.736699999999
}
foreach ($cell in $WVR_updates.Keys) {
  $ws.Range($cell).Value = $WVR_updates[$cell]
}

# Cells removed entirely in the target revision (no longer populated)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N96").ClearContents()

Write-Output "edits applied"